# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083   # subscript-three character used in PEPE price (0.0<sub>3</sub>0923)

$ws.Range("D2").Value = "61.810.18"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "2.400.80"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("D5").Value = "'560.51"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'141.89"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("E13").Value = "  -3.40%  "
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "2.830.13"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").Value = "61.803.06"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "2.398.86"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "'11.15"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'320.61"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.79"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "'4.10"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'65.44"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "'1.71"
$ws.Range("E24").Value = "  -3.80%  "
$ws.Range("D25").Value = "'8.65"
$ws.Range("E25").Value = "  -5.19%  "
$ws.Range("D26").Value = "'562.69"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").Value = "'0.0" + $sub3 + "0923"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "'8.12"
$ws.Range("E30").Value = "  -3.37%  "
$ws.Range("E31").Value = "  -4.96%  "
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").Value = "'1.87"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").Value = "'1.49"
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").Value = "'151.99"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("D38").Value = "'5.41"
$ws.Range("E38").Value = "  -5.28%  "
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("E41").Value = "  -5.72%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'147.17"
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("D44").Value = "'2.22"
$ws.Range("E44").Value = "  -4.19%  "
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("E46").Value = "  -3.01%  "
$ws.Range("D47").Value = "'19.74"
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").Value = "'0.0916"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("E51").Value = "  +0.28%  "
